$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values - repulled data
$updates = @{
    2  = -5
    3  = 4
    4  = 1
    5  = 3
    6  = 6
    7  = -1
    8  = -1
    9  = 7
    10 = -2
    12 = 5
    13 = -2
    14 = -5
    16 = -6
    18 = -3
    19 = -6
    20 = 2
    21 = -2
    22 = -1
    23 = -1
    24 = 2
    25 = 2
    26 = 5
    27 = 3
    29 = 9
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
